$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 174, shifting existing row 174 (and everything
# below it) down by one — matches the dimension growing from R267 to R268.
$ws.Rows("174:174").Insert(-4121)

# Populate the newly inserted row 174 with the new record.
$ws.Range("A174").Value = 3
$ws.Range("B174").Value = "Femacal de La Calera"
$ws.Range("C174").Value = "Coquimbo"
$ws.Range("D174").Value = 44572
$ws.Range("E174").Value = 5
$ws.Range("F174").Value = 100112009
$ws.Range("G174").Value = "Acelga"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 310
$ws.Range("K174").Value = 2300
$ws.Range("L174").Value = 2500
$ws.Range("M174").Value = 2397
$ws.Range("N174").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O174").Value = "Provincia de Quillota"
$ws.Range("P174").Value = 400
$ws.Range("Q174").Value = 6
$ws.Range("R174").Value = "Hortaliza"

# Match the date-formatted style used by the other rows' Fecha column.
$ws.Range("D174").NumberFormat = $ws.Range("D175").NumberFormat
